$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: LIKING_PROMPT - add trailing period to both languages ---
$ws.Range("B18").Value = "Das Stück hat mir gefallen."
$ws.Range("C18").Value = "I liked the piece."

# --- Row 19: DIFFICULTY_PROMPT - add trailing period to German text ---
$ws.Range("B19").Value = "Die Aufgabe war schwierig."

# --- Row 22: CONTINUE_MAIN_TEST - new button text ---
$ws.Range("B22").Value = "Beginne mit dem Experiment."
$ws.Range("C22").Value = "Start to the experiment."

# --- New rows 23-28: LIKERT1..LIKERT6 scale labels ---
$ws.Range("A23").Value = "LIKERT1"
$ws.Range("B23").Value = "Trifft gar nicht zu "
$ws.Range("C23").Value = "Completely disagree"

$ws.Range("A24").Value = "LIKERT2"
$ws.Range("B24").Value = "Trifft nicht zu"
$ws.Range("C24").Value = "Strongly disagree"

$ws.Range("A25").Value = "LIKERT3"
$ws.Range("B25").Value = "Trifft eher nicht zu"
$ws.Range("C25").Value = "Disagree"

$ws.Range("A26").Value = "LIKERT4"
$ws.Range("B26").Value = "Trifft eher zu"
$ws.Range("C26").Value = "Agree"

$ws.Range("A27").Value = "LIKERT5"
$ws.Range("B27").Value = "Trifft zu"
$ws.Range("C27").Value = "Strongly agree"

$ws.Range("A28").Value = "LIKERT6"
$ws.Range("B28").Value = "Triff sehr zu"
$ws.Range("C28").Value = "Completely agree"

# Style column A for the new rows the same as the rest of column A (style index 2:
# vertical top alignment) - copy format down from the row above.
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23:A28").PasteSpecial(-4122) | Out-Null

# Column B (D-G too) of rows 23-28: horizontal-left alignment, default font/size.
$ws.Range("B23:G28").HorizontalAlignment = -4131

# Column C of rows 23-28: horizontal-left alignment + bigger, black font (matches
# the LIKERT header emphasis used in the source workbook).
$ws.Range("C23:C28").Font.Size = 12
$ws.Range("C23:C28").Font.Color = 0

# Row heights for the new Likert rows.
$ws.Rows("23:28").RowHeight = 15.75

$excel.CutCopyMode = 0

# --- Update the view / selection to match the edited workbook state ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C19").Select()
